$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3442397.08
$ws.Range("C9").Value = 537276.99
$ws.Range("D9").Value = 3979674.07
$ws.Range("E9").Value = 13.50052744394718
$ws.Range("F9").Value = 86.49947255605281
$ws.Range("G9").Value = -48.07485790307787
$ws.Range("H9").Value = -37.83506695281292
$ws.Range("I9").Value = 34494
$ws.Range("J9").Value = 1462
$ws.Range("K9").Value = 35956
$ws.Range("L9").Value = 24812
$ws.Range("M9").Value = 160.3931190552958
$ws.Range("N9").Value = 9.503649222698595
